$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.09406588541035354
$ws.Range("D2").Value = 0.7813990817587921
$ws.Range("E2").Value = 0.03103577933686807
$ws.Range("F2").Value = 21.85007631241808
$ws.Range("G2").Value = 0.003150467317747274
$ws.Range("J2").Value = 0.2498563911488247
$ws.Range("L2").Value = 0.1341044033651713
$ws.Range("M2").Value = 10.42160049316291
$ws.Range("N2").Value = 1.504202102293192
$ws.Range("C3").Value = 0.09504084073810048
$ws.Range("D3").Value = 0.7628433209850698
$ws.Range("E3").Value = 0.0269115561100719
$ws.Range("F3").Value = 21.80513693331193
$ws.Range("G3").Value = 0.003174326715146544
$ws.Range("J3").Value = 0.2543071639378773
$ws.Range("L3").Value = 0.1220472746867642
$ws.Range("M3").Value = 10.18886034099157
$ws.Range("N3").Value = 1.490789305639211
$ws.Range("C4").Value = 0.09566797195529197
$ws.Range("D4").Value = 0.7522647551483601
$ws.Range("E4").Value = 0.02439388430628142
$ws.Range("F4").Value = 21.79565986186299
$ws.Range("G4").Value = 0.003189637043959326
$ws.Range("J4").Value = 0.2571821158558087
$ws.Range("L4").Value = 0.1147621229389415
$ws.Range("M4").Value = 10.05286381168023
$ws.Range("N4").Value = 1.483465168856085
$ws.Range("C5").Value = 0.09593073196163182
$ws.Range("D5").Value = 0.748152812020237
$ws.Range("E5").Value = 0.02337107268394334
$ws.Range("F5").Value = 21.79627632537057
$ws.Range("G5").Value = 0.003196043743431468
$ws.Range("J5").Value = 0.2583893534272015
$ws.Range("L5").Value = 0.1118212780225463
$ws.Range("M5").Value = 9.999146149881028
$ws.Range("N5").Value = 1.480711048301671
$ws.Range("C6").Value = 0.09597479875996129
$ws.Range("D6").Value = 0.74748187085072
$ws.Range("E6").Value = 0.02320141177935753
$ws.Range("F6").Value = 21.79664699046265
$ws.Range("G6").Value = 0.003197117737774999
$ws.Range("J6").Value = 0.2585919666402461
$ws.Range("L6").Value = 0.1113345890559572
$ws.Range("M6").Value = 9.990328188206377
$ws.Range("N6").Value = 1.480267692665919
$ws.Range("C7").Value = 0.09567148643887613
$ws.Range("D7").Value = 0.7522085023741738
$ws.Range("E7").Value = 0.02438007818238219
$ws.Range("F7").Value = 21.79565014175415
$ws.Range("G7").Value = 0.003189722766523894
$ws.Range("J7").Value = 0.2571982527477239
$ws.Range("L7").Value = 0.1147223509172051
$ws.Range("M7").Value = 10.05213250970226
$ws.Range("N7").Value = 1.483427090594333
$ws.Range("C8").Value = 0.09439615292400561
$ws.Range("D8").Value = 0.7748282793598946
$ws.Range("E8").Value = 0.02961039513922969
$ws.Range("F8").Value = 21.83077382663072
$ws.Range("G8").Value = 0.003158557861903776
$ws.Range("J8").Value = 0.2513614658790555
$ws.Range("L8").Value = 0.1299215390404953
$ws.Range("M8").Value = 10.33989661738667
$ws.Range("N8").Value = 1.49938899759384
$ws.Range("C9").Value = 0.09211992855199469
$ws.Range("D9").Value = 0.8259347307075586
$ws.Range("E9").Value = 0.04000776789959559
$ws.Range("F9").Value = 22.04706602147974
$ws.Range("G9").Value = 0.003102614185249832
$ws.Range("J9").Value = 0.2410469725006621
$ws.Range("L9").Value = 0.1607472741385294
$ws.Range("M9").Value = 10.96069862195753
$ws.Range("N9").Value = 1.537868943929169
$ws.Range("C10").Value = 0.09058243514575359
$ws.Range("D10").Value = 0.8680159266000942
$ws.Range("E10").Value = 0.0477691647398828
$ws.Range("F10").Value = 22.30125941166455
$ws.Range("G10").Value = 0.003064562736744471
$ws.Range("J10").Value = 0.2341635555831978
$ws.Range("L10").Value = 0.1841414497032474
$ws.Range("M10").Value = 11.45379173666947
$ws.Range("N10").Value = 1.570455656400355
$ws.Range("C11").Value = 0.08991179137747274
$ws.Range("D11").Value = 0.8882423224589502
$ws.Range("E11").Value = 0.05133520238123879
$ws.Range("F11").Value = 22.43886286086257
$ws.Range("G11").Value = 0.003047891440407402
$ws.Range("J11").Value = 0.2311841363844493
$ws.Range("L11").Value = 0.1949754026417878
$ws.Range("M11").Value = 11.68674976435551
$ws.Range("N11").Value = 1.586206036710678
$ws.Range("C12").Value = 0.08966193432079628
$ws.Range("D12").Value = 0.8960659223540688
$ws.Range("E12").Value = 0.05269139294244241
$ws.Range("F12").Value = 22.49424104037342
$ws.Range("G12").Value = 0.003041668332079648
$ws.Range("J12").Value = 0.2300778678330815
$ws.Range("L12").Value = 0.1991080599313193
$ws.Range("M12").Value = 11.77626110260366
$ws.Range("N12").Value = 1.592302574001849
$ws.Range("C13").Value = 0.08971556369055023
$ws.Range("D13").Value = 0.8943735218377924
$ws.Range("E13").Value = 0.05239904260448469
$ws.Range("F13").Value = 22.48216704183102
$ws.Range("G13").Value = 0.00304300461993773
$ws.Range("J13").Value = 0.2303151431002277
$ws.Range("L13").Value = 0.198216639401366
$ws.Range("M13").Value = 11.75692477040104
$ws.Range("N13").Value = 1.590983711120714
$ws.Range("C14").Value = 0.08989115350239452
$ws.Range("D14").Value = 0.8888826313265099
$ws.Range("E14").Value = 0.05144665623990363
$ws.Range("F14").Value = 22.44335264312247
$ws.Range("G14").Value = 0.00304737767010721
$ws.Range("J14").Value = 0.2310926818434389
$ws.Range("L14").Value = 0.195314781305413
$ws.Range("M14").Value = 11.69408765012213
$ws.Range("N14").Value = 1.586704958438588
$ws.Range("C15").Value = 0.08999924034207751
$ws.Range("D15").Value = 0.8855409604997817
$ws.Range("E15").Value = 0.05086407100142765
$ws.Range("F15").Value = 22.42000707757012
$ws.Range("G15").Value = 0.003050067946904989
$ws.Range("J15").Value = 0.2315718123817767
$ws.Range("L15").Value = 0.1935413024905017
$ws.Range("M15").Value = 11.65576836434661
$ws.Range("N15").Value = 1.584101287309636
$ws.Range("C16").Value = 0.09062683917499292
$ws.Range("D16").Value = 0.8667166013050291
$ws.Range("E16").Value = 0.04753688814022894
$ws.Range("F16").Value = 22.29271816723895
$ws.Range("G16").Value = 0.003065664935342776
$ws.Range("J16").Value = 0.2343613349885239
$ws.Range("L16").Value = 0.1834374948883664
$ws.Range("M16").Value = 11.43874573239157
$ws.Range("N16").Value = 1.569444887075775
$ws.Range("C17").Value = 0.09101919420266569
$ws.Range("D17").Value = 0.8554518570684309
$ws.Range("E17").Value = 0.04550532480422476
$ws.Range("F17").Value = 22.22033472353371
$ws.Range("G17").Value = 0.003075395413924471
$ws.Range("J17").Value = 0.2361116085620338
$ws.Range("L17").Value = 0.1772899924172719
$ws.Range("M17").Value = 11.30785991992502
$ws.Range("N17").Value = 1.560690216217921
$ws.Range("C18").Value = 0.09124757659378169
$ws.Range("D18").Value = 0.849074115155247
$ws.Range("E18").Value = 0.04434010354536611
$ws.Range("F18").Value = 22.18076249339316
$ws.Range("G18").Value = 0.003081052342298456
$ws.Range("J18").Value = 0.2371326069733612
$ws.Range("L18").Value = 0.1737720093636312
$ws.Range("M18").Value = 11.23338851311442
$ws.Range("N18").Value = 1.555742026665428
$ws.Range("C19").Value = 0.09132536949311643
$ws.Range("D19").Value = 0.8469318929663814
$ws.Range("E19").Value = 0.04394612155578415
$ws.Range("F19").Value = 22.16771467407227
$ws.Range("G19").Value = 0.00308297807861221
$ws.Range("J19").Value = 0.2374807494760702
$ws.Range("L19").Value = 0.1725838812464957
$ws.Range("M19").Value = 11.20831148466249
$ws.Range("N19").Value = 1.554081674878432
$ws.Range("C20").Value = 0.09097714710454774
$ws.Range("D20").Value = 0.8566404495642814
$ws.Range("E20").Value = 0.04572124356261043
$ws.Range("F20").Value = 22.22782599213548
$ws.Range("G20").Value = 0.003074353369131722
$ws.Range("J20").Value = 0.2359238093373826
$ws.Range("L20").Value = 0.1779425346345107
$ws.Range("M20").Value = 11.32170867526466
$ws.Range("N20").Value = 1.561613141838052
$ws.Range("C21").Value = 0.08983946750527849
$ws.Range("D21").Value = 0.8904909095277844
$ws.Range("E21").Value = 0.0517262315161986
$ws.Range("F21").Value = 22.45466368287651
$ws.Range("G21").Value = 0.003046090773230426
$ws.Range("J21").Value = 0.2308637024273885
$ws.Range("L21").Value = 0.1961662905619193
$ws.Range("M21").Value = 11.71250885537739
$ws.Range("N21").Value = 1.587958151499748
$ws.Range("C22").Value = 0.08911981541861458
$ws.Range("D22").Value = 0.9135761440163606
$ws.Range("E22").Value = 0.05568507011862778
$ws.Range("F22").Value = 22.62203638740476
$ws.Range("G22").Value = 0.003028142904848781
$ws.Range("J22").Value = 0.2276847730303011
$ws.Range("L22").Value = 0.2082531314368339
$ws.Range("M22").Value = 11.9754952274248
$ws.Range("N22").Value = 1.60594633188515
$ws.Range("C23").Value = 0.08950173376541848
$ws.Range("D23").Value = 0.9011641517977864
$ws.Range("E23").Value = 0.05356877924592141
$ws.Range("F23").Value = 22.53091816706717
$ws.Range("G23").Value = 0.003037674767958247
$ws.Range("J23").Value = 0.2293696560898333
$ws.Range("L23").Value = 0.2017851387559517
$ws.Range("M23").Value = 11.8344232682183
$ws.Range("N23").Value = 1.596275549543321
$ws.Range("C24").Value = 0.09099614783520593
$ws.Range("D24").Value = 0.8561027806985635
$ws.Range("E24").Value = 0.04562361817588112
$ws.Range("F24").Value = 22.22443283765648
$ws.Range("G24").Value = 0.003074824282173135
$ws.Range("J24").Value = 0.2360086674558204
$ws.Range("L24").Value = 0.1776474697202985
$ws.Range("M24").Value = 11.31544523793761
$ws.Range("N24").Value = 1.561195622370235
$ws.Range("C25").Value = 0.09271186910827467
$ws.Range("D25").Value = 0.8113446580392178
$ws.Range("E25").Value = 0.03717626822780318
$ws.Range("F25").Value = 21.97223075996675
$ws.Range("G25").Value = 0.003117205193618258
$ws.Range("J25").Value = 0.2437157102298162
$ws.Range("L25").Value = 0.1522864191493341
$ws.Range("M25").Value = 10.78645306983563
$ws.Range("N25").Value = 1.504202102293192
